$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.690.85'
$ws.Range("E2").Value = '  -6.42%  '
$ws.Range("D3").Value = '2.895.73'
$ws.Range("E3").Value = '  -4.63%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.84'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -5.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '121.68'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -6.08%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '2.887.21'
$ws.Range("E8").Value = '  -4.86%  '
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("E10").Value = '  -8.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.75'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -9.36%  '
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("E13").Value = '  -8.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.62'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.39%  '
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").Value = '3.370.06'
$ws.Range("E16").Value = '  -4.50%  '
$ws.Range("D17").Value = '2.892.30'
$ws.Range("E17").Value = '  -4.68%  '
$ws.Range("D18").Value = '57.633.91'
$ws.Range("E18").Value = '  -6.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.46'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '406.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -8.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.78'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.653'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.32%  '
$ws.Range("E23").Value = '  -7.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.59'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '76.81'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.46%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("E28").Value = '  -3.66%  '
$ws.Range("E29").Value = '  -3.49%  '
$ws.Range("E30").Value = '  -4.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.01'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '24.57'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0948'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("E34").Value = '  -12.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.899'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.34'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '48.38'
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.40'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +6.69%  '
$ws.Range("D39").Value = '0.0₃0616'
$ws.Range("E39").Value = '  -11.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0342'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.73%  '
$ws.Range("E41").Value = '  -4.55%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.610.49'
$ws.Range("E42").Value = '  -2.46%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '358.88'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.34'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -7.69%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '117.66'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.227'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.45%  '
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.93'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.39'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.94'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.28%  '
